# EPBDS-12620 Difference in error response structure between kafka and webservice call
#
# The "Good Night" label in cell E11 of the "Rules" sheet is replaced with the
# literal text  = error("fail")  (entered as text, not as a formula, hence the
# leading apostrophe which Excel stores as quotePrefix="1" on the cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E11").Value = "'= error(""fail"")"

# Cosmetic: leave the selection where the author left it when they saved.
$null = $ws.Range("F9").Select()
